# Update FNO reversal scanner plan: refresh tickers, prices, and position% for rows 2-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @(2, "BAJAJ-AUTO", 9113.405025202483, 9.30232558139535),
    @(3, "TI", 467.5684548757732, 9.30232558139535),
    @(4, "PVRINOX", 1117.750226124427, 6.976744186046512),
    @(5, "LTF", 224.0787602742681, 6.976744186046512),
    @(6, "PIIND", 4595.026806698306, 2.325581395348837),
    @(7, "MASTEK", 2998.999843288399, 2.325581395348837),
    @(8, "SCHAEFFLER", 4747.486630466043, 2.325581395348837),
    @(9, "UPL", 746.5463222727701, 6.976744186046512),
    @(10, "CRAFTSMAN", 7111.741741431812, 2.325581395348837),
    @(11, "ASHAPURMIN", 622.7105881821097, 4.651162790697675),
    @(12, "CHEMPLASTS", 479.8745961253245, 6.976744186046512),
    @(13, "METROBRAND", 1300.862736668633, 2.325581395348837),
    @(14, "SUNTV", 634.4313398645224, 2.325581395348837),
    @(15, "POONAWALLA", 490.3154832072964, 2.325581395348837),
    @(16, "TATAINVEST", 7184.423831922282, 4.651162790697675),
    @(17, "TIMKEN", 3561.631160176454, 6.976744186046512),
    @(18, "PRAKASH", 185.7201128938527, 6.976744186046512),
    @(19, "TATASTEEL", 168.7678033002403, 4.651162790697675),
    @(20, "FLUOROCHEM", 3653.500678453571, 4.651162790697675),
    @(21, "VBL", 499.5081586398466, 4.651162790697675)
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Cells.Item($row, 2).Value = $u[1]
    $ws.Cells.Item($row, 3).Value = $u[2]
    $ws.Cells.Item($row, 4).Value = $u[3]
}
